$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row captions: "<...>_old" -> "<...>_FV2404" and
#    "<...>_new" -> "<...>_FV2410" (the "diff" header is left untouched).
$headerRange = $ws.Range("A1:U1")
for ($i = 1; $i -le $headerRange.Columns.Count; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $text = [string]$cell.Value()
    if ($text.EndsWith("_old")) {
        $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2404"
    } elseif ($text.EndsWith("_new")) {
        $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2410"
    }
}

# 2. Turn the used range into a real Excel table ("Table1") so the header
#    row gets filter buttons and the structured-table styling.
$tableRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
